$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.026.22"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "2.634.17"
$ws.Range("E3").Value = "  +5.49%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'110.73"
$ws.Range("E5").Value = "  +4.78%  "
$ws.Range("D6").Value = "'320.99"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'0.520"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "'39.57"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").Value = "'19.89"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'0.0808"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "'7.22"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "3.045.98"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").Value = "2.657.91"
$ws.Range("E16").Value = "  +7.33%  "
$ws.Range("D17").Value = "'0.857"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").Value = "48.951.20"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").Value = "'12.88"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "'6.67"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").Value = "'2.89"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "0.0₃0942"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").Value = "'270.19"
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").Value = "'70.06"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'2.54"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").Value = "'26.27"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'10.09"
$ws.Range("E28").Value = "  +4.93%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").Value = "'35.45"
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("D31").Value = "'0.138"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'49.41"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "'5.42"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'19.17"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'0.0794"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("D37").Value = "'4.91"
$ws.Range("E37").Value = "  +8.80%  "
$ws.Range("D38").Value = "'2.02"
$ws.Range("E38").Value = "  +3.88%  "
$ws.Range("D39").Value = "'3.14"
$ws.Range("E39").Value = "  +8.65%  "
$ws.Range("D40").Value = "'124.79"
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'22.71"
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.111"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'0.0313"
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("D45").Value = "2.087.26"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("D46").Value = "'3.22"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'2.12"
$ws.Range("E47").Value = "  +8.62%  "
$ws.Range("E48").Value = "  +4.88%  "
$ws.Range("D49").Value = "2.900.11"
$ws.Range("E49").Value = "  +5.93%  "
$ws.Range("D50").Value = "'8.88"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "'59.01"
$ws.Range("E51").Value = "  +5.28%  "
